$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 7).Value = 5.573004
$ws.Cells.Item(2, 8).Value = 16.719012
$ws.Cells.Item(2, 9).Value = 0.1383950099922425
$ws.Cells.Item(2, 10).Value = 0.1494586813194652
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 13).Value = 3.946674666666667
$ws.Cells.Item(2, 14).Value = 11.840024
$ws.Cells.Item(2, 15).Value = 0.008082287850277156
$ws.Cells.Item(2, 16).Value = 0.008317683527585098
$ws.Cells.Item(2, 17).Value = 21.994833704032
$ws.Cells.Item(2, 18).Value = 197.953503336288
$ws.Cells.Item(2, 19).Value = 0.001118548307799287
$ws.Cells.Item(2, 20).Value = 0.001243150011665507

$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 7).Value = 5.573004
$ws.Cells.Item(3, 8).Value = 16.719012
$ws.Cells.Item(3, 9).Value = 0.1383950099922425
$ws.Cells.Item(3, 10).Value = 0.1494586813194652
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 13).Value = 185.8027443333333
$ws.Cells.Item(3, 14).Value = 557.408233
$ws.Cells.Item(3, 15).Value = 0.3805003933455167
$ws.Cells.Item(3, 16).Value = 0.3915824222792467
$ws.Cells.Item(3, 17).Value = 1035.479437380644
$ws.Cells.Item(3, 18).Value = 9319.314936425795
$ws.Cells.Item(3, 19).Value = 0.05265935573910498
$ws.Cells.Item(3, 20).Value = 0.0585253924617382

$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 7).Value = 5.573004
$ws.Cells.Item(4, 8).Value = 16.719012
$ws.Cells.Item(4, 9).Value = 0.1383950099922425
$ws.Cells.Item(4, 10).Value = 0.1494586813194652
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 13).Value = 117.3394243333333
$ws.Cells.Item(4, 14).Value = 352.018273
$ws.Cells.Item(4, 15).Value = 0.2402962199184265
$ws.Cells.Item(4, 16).Value = 0.2472948188906589
$ws.Cells.Item(4, 17).Value = 653.9330811673641
$ws.Cells.Item(4, 18).Value = 5885.397730506276
$ws.Cells.Item(4, 19).Value = 0.03325579775670873
$ws.Cells.Item(4, 20).Value = 0.03696035752853385

$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 7).Value = 5.573004
$ws.Cells.Item(5, 8).Value = 16.719012
$ws.Cells.Item(5, 9).Value = 0.1383950099922425
$ws.Cells.Item(5, 10).Value = 0.1494586813194652
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 13).Value = 139.7641523333333
$ws.Cells.Item(5, 14).Value = 419.292457
$ws.Cells.Item(5, 15).Value = 0.2862192112890951
$ws.Cells.Item(5, 16).Value = 0.2945553119511906
$ws.Cells.Item(5, 17).Value = 778.906180010276
$ws.Cells.Item(5, 18).Value = 7010.155620092484
$ws.Cells.Item(5, 19).Value = 0.03961131060632608
$ws.Cells.Item(5, 20).Value = 0.04402384849986866

$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 7).Value = 5.573004
$ws.Cells.Item(6, 8).Value = 16.719012
$ws.Cells.Item(6, 9).Value = 0.1383950099922425
$ws.Cells.Item(6, 10).Value = 0.1494586813194652
$ws.Cells.Item(6, 11).Value = 2
$ws.Cells.Item(6, 13).Value = 41.458574
$ws.Cells.Item(6, 14).Value = 82.917148
$ws.Cells.Item(6, 15).Value = 0.0849018875966847
$ws.Cells.Item(6, 16).Value = 0.05824976335131885
$ws.Cells.Item(6, 17).Value = 231.048798736296
$ws.Cells.Item(6, 18).Value = 1386.292792417776
$ws.Cells.Item(6, 19).Value = 0.01174999758230343
$ws.Cells.Item(6, 20).Value = 0.00870593281765903

$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 7).Value = 19.39389166666666
$ws.Cells.Item(7, 8).Value = 58.181675
$ws.Cells.Item(7, 9).Value = 0.4816106055184604
$ws.Cells.Item(7, 10).Value = 0.520111859627692
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 13).Value = 3.946674666666667
$ws.Cells.Item(7, 14).Value = 11.840024
$ws.Cells.Item(7, 15).Value = 0.008082287850277156
$ws.Cells.Item(7, 16).Value = 0.008317683527585098
$ws.Cells.Item(7, 17).Value = 76.5413809289111
$ws.Cells.Item(7, 18).Value = 688.8724283601999
$ws.Cells.Item(7, 19).Value = 0.003892515545546477
$ws.Cells.Item(7, 20).Value = 0.004326125847326907

$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 7).Value = 19.39389166666666
$ws.Cells.Item(8, 8).Value = 58.181675
$ws.Cells.Item(8, 9).Value = 0.4816106055184604
$ws.Cells.Item(8, 10).Value = 0.520111859627692
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(8, 13).Value = 185.8027443333333
$ws.Cells.Item(8, 14).Value = 557.408233
$ws.Cells.Item(8, 15).Value = 0.3805003933455167
$ws.Cells.Item(8, 16).Value = 0.3915824222792467
$ws.Cells.Item(8, 17).Value = 3603.43829497003
$ws.Cells.Item(8, 18).Value = 32430.94465473027
$ws.Cells.Item(8, 19).Value = 0.1832530248391467
$ws.Cells.Item(8, 20).Value = 0.2036666618491752

$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 7).Value = 19.39389166666666
$ws.Cells.Item(9, 8).Value = 58.181675
$ws.Cells.Item(9, 9).Value = 0.4816106055184604
$ws.Cells.Item(9, 10).Value = 0.520111859627692
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 13).Value = 117.3394243333333
$ws.Cells.Item(9, 14).Value = 352.018273
$ws.Cells.Item(9, 15).Value = 0.2402962199184265
$ws.Cells.Item(9, 16).Value = 0.2472948188906589
$ws.Cells.Item(9, 17).Value = 2275.668083749697
$ws.Cells.Item(9, 18).Value = 20481.01275374728
$ws.Cells.Item(9, 19).Value = 0.1157292079787105
$ws.Cells.Item(9, 20).Value = 0.1286209681295139

$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 7).Value = 19.39389166666666
$ws.Cells.Item(10, 8).Value = 58.181675
$ws.Cells.Item(10, 9).Value = 0.4816106055184604
$ws.Cells.Item(10, 10).Value = 0.520111859627692
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 13).Value = 139.7641523333333
$ws.Cells.Item(10, 14).Value = 419.292457
$ws.Cells.Item(10, 15).Value = 0.2862192112890951
$ws.Cells.Item(10, 16).Value = 0.2945553119511906
$ws.Cells.Item(10, 17).Value = 2710.570829236164
$ws.Cells.Item(10, 18).Value = 24395.13746312548
$ws.Cells.Item(10, 19).Value = 0.1378462076599572
$ws.Cells.Item(10, 20).Value = 0.1532017110621487

$ws.Cells.Item(11, 5).Value = 3
$ws.Cells.Item(11, 7).Value = 19.39389166666666
$ws.Cells.Item(11, 8).Value = 58.181675
$ws.Cells.Item(11, 9).Value = 0.4816106055184604
$ws.Cells.Item(11, 10).Value = 0.520111859627692
$ws.Cells.Item(11, 11).Value = 2
$ws.Cells.Item(11, 13).Value = 41.458574
$ws.Cells.Item(11, 14).Value = 82.917148
$ws.Cells.Item(11, 15).Value = 0.0849018875966847
$ws.Cells.Item(11, 16).Value = 0.05824976335131885
$ws.Cells.Item(11, 17).Value = 804.0430928104832
$ws.Cells.Item(11, 18).Value = 4824.2585568629
$ws.Cells.Item(11, 19).Value = 0.04088964949509958
$ws.Cells.Item(11, 20).Value = 0.03029639273952743

$ws.Cells.Item(12, 5).Value = 3
$ws.Cells.Item(12, 7).Value = 3.432072333333333
$ws.Cells.Item(12, 8).Value = 10.296217
$ws.Cells.Item(12, 9).Value = 0.08522902277941405
$ws.Cells.Item(12, 10).Value = 0.09204246132481156
$ws.Cells.Item(12, 11).Value = 3
$ws.Cells.Item(12, 13).Value = 3.946674666666667
$ws.Cells.Item(12, 14).Value = 11.840024
$ws.Cells.Item(12, 15).Value = 0.008082287850277156
$ws.Cells.Item(12, 16).Value = 0.008317683527585098
$ws.Cells.Item(12, 17).Value = 13.54527293213422
$ws.Cells.Item(12, 18).Value = 121.907456389208
$ws.Cells.Item(12, 19).Value = 0.0006888454953010532
$ws.Cells.Item(12, 20).Value = 0.0007655800643997736

$ws.Cells.Item(13, 5).Value = 3
$ws.Cells.Item(13, 7).Value = 3.432072333333333
$ws.Cells.Item(13, 8).Value = 10.296217
$ws.Cells.Item(13, 9).Value = 0.08522902277941405
$ws.Cells.Item(13, 10).Value = 0.09204246132481156
$ws.Cells.Item(13, 11).Value = 3
$ws.Cells.Item(13, 13).Value = 185.8027443333333
$ws.Cells.Item(13, 14).Value = 557.408233
$ws.Cells.Item(13, 15).Value = 0.3805003933455167
$ws.Cells.Item(13, 16).Value = 0.3915824222792467
$ws.Cells.Item(13, 17).Value = 637.68845828384
$ws.Cells.Item(13, 18).Value = 5739.19612455456
$ws.Cells.Item(13, 19).Value = 0.03242967669202104
$ws.Cells.Item(13, 20).Value = 0.03604220995811359

$ws.Cells.Item(14, 5).Value = 3
$ws.Cells.Item(14, 7).Value = 3.432072333333333
$ws.Cells.Item(14, 8).Value = 10.296217
$ws.Cells.Item(14, 9).Value = 0.08522902277941405
$ws.Cells.Item(14, 10).Value = 0.09204246132481156
$ws.Cells.Item(14, 11).Value = 3
$ws.Cells.Item(14, 13).Value = 117.3394243333333
$ws.Cells.Item(14, 14).Value = 352.018273
$ws.Cells.Item(14, 15).Value = 0.2402962199184265
$ws.Cells.Item(14, 16).Value = 0.2472948188906589
$ws.Cells.Item(14, 17).Value = 402.7173918636934
$ws.Cells.Item(14, 18).Value = 3624.456526773241
$ws.Cells.Item(14, 19).Value = 0.02048021200123466
$ws.Cells.Item(14, 20).Value = 0.02276162380356975

$ws.Cells.Item(15, 5).Value = 3
$ws.Cells.Item(15, 7).Value = 3.432072333333333
$ws.Cells.Item(15, 8).Value = 10.296217
$ws.Cells.Item(15, 9).Value = 0.08522902277941405
$ws.Cells.Item(15, 10).Value = 0.09204246132481156
$ws.Cells.Item(15, 11).Value = 3
$ws.Cells.Item(15, 13).Value = 139.7641523333333
$ws.Cells.Item(15, 14).Value = 419.292457
$ws.Cells.Item(15, 15).Value = 0.2862192112890951
$ws.Cells.Item(15, 16).Value = 0.2945553119511906
$ws.Cells.Item(15, 17).Value = 479.6806804150187
$ws.Cells.Item(15, 18).Value = 4317.126123735168
$ws.Cells.Item(15, 19).Value = 0.0243941836788642
$ws.Cells.Item(15, 20).Value = 0.02711159590828526

$ws.Cells.Item(16, 5).Value = 3
$ws.Cells.Item(16, 7).Value = 3.432072333333333
$ws.Cells.Item(16, 8).Value = 10.296217
$ws.Cells.Item(16, 9).Value = 0.08522902277941405
$ws.Cells.Item(16, 10).Value = 0.09204246132481156
$ws.Cells.Item(16, 11).Value = 2
$ws.Cells.Item(16, 13).Value = 41.458574
$ws.Cells.Item(16, 14).Value = 82.917148
$ws.Cells.Item(16, 15).Value = 0.0849018875966847
$ws.Cells.Item(16, 16).Value = 0.05824976335131885
$ws.Cells.Item(16, 17).Value = 142.2888248048526
$ws.Cells.Item(16, 18).Value = 853.7329488291159
$ws.Cells.Item(16, 19).Value = 0.007236104911993091
$ws.Cells.Item(16, 20).Value = 0.005361451590443191

$ws.Cells.Item(17, 5).Value = 3
$ws.Cells.Item(17, 7).Value = 2.927161333333333
$ws.Cells.Item(17, 8).Value = 8.781483999999999
$ws.Cells.Item(17, 9).Value = 0.0726905134063375
$ws.Cells.Item(17, 10).Value = 0.07850158960756669
$ws.Cells.Item(17, 11).Value = 3
$ws.Cells.Item(17, 13).Value = 3.946674666666667
$ws.Cells.Item(17, 14).Value = 11.840024
$ws.Cells.Item(17, 15).Value = 0.008082287850277156
$ws.Cells.Item(17, 16).Value = 0.008317683527585098
$ws.Cells.Item(17, 17).Value = 11.55255347951289
$ws.Cells.Item(17, 18).Value = 103.972981315616
$ws.Cells.Item(17, 19).Value = 0.0005875056533344503
$ws.Cells.Item(17, 20).Value = 0.000652951378768103

$ws.Cells.Item(18, 5).Value = 3
$ws.Cells.Item(18, 7).Value = 2.927161333333333
$ws.Cells.Item(18, 8).Value = 8.781483999999999
$ws.Cells.Item(18, 9).Value = 0.0726905134063375
$ws.Cells.Item(18, 10).Value = 0.07850158960756669
$ws.Cells.Item(18, 11).Value = 3
$ws.Cells.Item(18, 13).Value = 185.8027443333333
$ws.Cells.Item(18, 14).Value = 557.408233
$ws.Cells.Item(18, 15).Value = 0.3805003933455167
$ws.Cells.Item(18, 16).Value = 0.3915824222792467
$ws.Cells.Item(18, 17).Value = 543.8746088397523
$ws.Cells.Item(18, 18).Value = 4894.871479557772
$ws.Cells.Item(18, 19).Value = 0.02765876894359897
$ws.Cells.Item(18, 20).Value = 0.0307398426113023

$ws.Cells.Item(19, 5).Value = 3
$ws.Cells.Item(19, 7).Value = 2.927161333333333
$ws.Cells.Item(19, 8).Value = 8.781483999999999
$ws.Cells.Item(19, 9).Value = 0.0726905134063375
$ws.Cells.Item(19, 10).Value = 0.07850158960756669
$ws.Cells.Item(19, 11).Value = 3
$ws.Cells.Item(19, 13).Value = 117.3394243333333
$ws.Cells.Item(19, 14).Value = 352.018273
$ws.Cells.Item(19, 15).Value = 0.2402962199184265
$ws.Cells.Item(19, 16).Value = 0.2472948188906589
$ws.Cells.Item(19, 17).Value = 343.4714257841258
$ws.Cells.Item(19, 18).Value = 3091.242832057132
$ws.Cells.Item(19, 19).Value = 0.0174672555954726
$ws.Cells.Item(19, 20).Value = 0.01941303638463203

$ws.Cells.Item(20, 5).Value = 3
$ws.Cells.Item(20, 7).Value = 2.927161333333333
$ws.Cells.Item(20, 8).Value = 8.781483999999999
$ws.Cells.Item(20, 9).Value = 0.0726905134063375
$ws.Cells.Item(20, 10).Value = 0.07850158960756669
$ws.Cells.Item(20, 11).Value = 3
$ws.Cells.Item(20, 13).Value = 139.7641523333333
$ws.Cells.Item(20, 14).Value = 419.292457
$ws.Cells.Item(20, 15).Value = 0.2862192112890951
$ws.Cells.Item(20, 16).Value = 0.2945553119511906
$ws.Cells.Item(20, 17).Value = 409.112222496243
$ws.Cells.Item(20, 18).Value = 3682.010002466188
$ws.Cells.Item(20, 19).Value = 0.02080542141536131
$ws.Cells.Item(20, 20).Value = 0.02312306021552114

$ws.Cells.Item(21, 5).Value = 3
$ws.Cells.Item(21, 7).Value = 2.927161333333333
$ws.Cells.Item(21, 8).Value = 8.781483999999999
$ws.Cells.Item(21, 9).Value = 0.0726905134063375
$ws.Cells.Item(21, 10).Value = 0.07850158960756669
$ws.Cells.Item(21, 11).Value = 2
$ws.Cells.Item(21, 13).Value = 41.458574
$ws.Cells.Item(21, 14).Value = 82.917148
$ws.Cells.Item(21, 15).Value = 0.0849018875966847
$ws.Cells.Item(21, 16).Value = 0.05824976335131885
$ws.Cells.Item(21, 17).Value = 121.3559347479386
$ws.Cells.Item(21, 18).Value = 728.1356084876319
$ws.Cells.Item(21, 19).Value = 0.006171561798570169
$ws.Cells.Item(21, 20).Value = 0.004572699017343111

$ws.Cells.Item(22, 5).Value = 2
$ws.Cells.Item(22, 7).Value = 8.9426925
$ws.Cells.Item(22, 8).Value = 17.885385
$ws.Cells.Item(22, 9).Value = 0.2220748483035455
$ws.Cells.Item(22, 10).Value = 0.1598854081204645
$ws.Cells.Item(22, 11).Value = 3
$ws.Cells.Item(22, 13).Value = 3.946674666666667
$ws.Cells.Item(22, 14).Value = 11.840024
$ws.Cells.Item(22, 15).Value = 0.008082287850277156
$ws.Cells.Item(22, 16).Value = 0.008317683527585098
$ws.Cells.Item(22, 17).Value = 35.29389794154
$ws.Cells.Item(22, 18).Value = 211.76338764924
$ws.Cells.Item(22, 19).Value = 0.001794872848295888
$ws.Cells.Item(22, 20).Value = 0.001329876225424809

$ws.Cells.Item(23, 5).Value = 2
$ws.Cells.Item(23, 7).Value = 8.9426925
$ws.Cells.Item(23, 8).Value = 17.885385
$ws.Cells.Item(23, 9).Value = 0.2220748483035455
$ws.Cells.Item(23, 10).Value = 0.1598854081204645
$ws.Cells.Item(23, 11).Value = 3
$ws.Cells.Item(23, 13).Value = 185.8027443333333
$ws.Cells.Item(23, 14).Value = 557.408233
$ws.Cells.Item(23, 15).Value = 0.3805003933455167
$ws.Cells.Item(23, 16).Value = 0.3915824222792467
$ws.Cells.Item(23, 17).Value = 1661.576808229117
$ws.Cells.Item(23, 18).Value = 9969.460849374706
$ws.Cells.Item(23, 19).Value = 0.08449956713164498
$ws.Cells.Item(23, 20).Value = 0.06260831539891745

$ws.Cells.Item(24, 5).Value = 2
$ws.Cells.Item(24, 7).Value = 8.9426925
$ws.Cells.Item(24, 8).Value = 17.885385
$ws.Cells.Item(24, 9).Value = 0.2220748483035455
$ws.Cells.Item(24, 10).Value = 0.1598854081204645
$ws.Cells.Item(24, 11).Value = 3
$ws.Cells.Item(24, 13).Value = 117.3394243333333
$ws.Cells.Item(24, 14).Value = 352.018273
$ws.Cells.Item(24, 15).Value = 0.2402962199184265
$ws.Cells.Item(24, 16).Value = 0.2472948188906589
$ws.Cells.Item(24, 17).Value = 1049.330389940018
$ws.Cells.Item(24, 18).Value = 6295.982339640105
$ws.Cells.Item(24, 19).Value = 0.05336374658629995
$ws.Cells.Item(24, 20).Value = 0.03953883304440935

$ws.Cells.Item(25, 5).Value = 2
$ws.Cells.Item(25, 7).Value = 8.9426925
$ws.Cells.Item(25, 8).Value = 17.885385
$ws.Cells.Item(25, 9).Value = 0.2220748483035455
$ws.Cells.Item(25, 10).Value = 0.1598854081204645
$ws.Cells.Item(25, 11).Value = 3
$ws.Cells.Item(25, 13).Value = 139.7641523333333
$ws.Cells.Item(25, 14).Value = 419.292457
$ws.Cells.Item(25, 15).Value = 0.2862192112890951
$ws.Cells.Item(25, 16).Value = 0.2945553119511906
$ws.Cells.Item(25, 17).Value = 1249.867836840157
$ws.Cells.Item(25, 18).Value = 7499.207021040945
$ws.Cells.Item(25, 19).Value = 0.06356208792858621
$ws.Cells.Item(25, 20).Value = 0.04709509626536685

$ws.Cells.Item(26, 5).Value = 2
$ws.Cells.Item(26, 7).Value = 8.9426925
$ws.Cells.Item(26, 8).Value = 17.885385
$ws.Cells.Item(26, 9).Value = 0.2220748483035455
$ws.Cells.Item(26, 10).Value = 0.1598854081204645
$ws.Cells.Item(26, 11).Value = 2
$ws.Cells.Item(26, 13).Value = 41.458574
$ws.Cells.Item(26, 14).Value = 82.917148
$ws.Cells.Item(26, 15).Value = 0.0849018875966847
$ws.Cells.Item(26, 16).Value = 0.05824976335131885
$ws.Cells.Item(26, 17).Value = 370.751278770495
$ws.Cells.Item(26, 18).Value = 1483.00511508198
$ws.Cells.Item(26, 19).Value = 0.01885457380871842
$ws.Cells.Item(26, 20).Value = 0.009313287186346093
